$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename existing row 4 label (Biased variant)
$ws.Range("A4").Value = "Cropped_Compensated_NaivePartition_ Biased"

# Add new row 5: Unbiased variant
$ws.Range("A5").Value = "Cropped_Compensated_NaivePartition_ Unbiased"
$ws.Range("B5").Value = 0.2345
$ws.Range("C5").Value = 0.2722
$ws.Range("D5").Value = 0.3338
$ws.Range("E5").Value = 0.3854

# Add new row 6: 5Fold partition
$ws.Range("A6").Value = "Cropped_Compensated_5FoldPartition"
$ws.Range("B6").Value = 0.3539
$ws.Range("C6").Value = 0.36
$ws.Range("D6").Value = 0.435
$ws.Range("E6").Value = 0.7188

# Update the selected cell to match the final saved selection
$ws.Range("J9").Select()
